# Croatia HNL base update — 29-03-2024 13:24
# Appends 4 new fixtures (rows 136-139) below the existing data (which
# ends at row 135), matching the layout/format of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy the per-column formatting of the last existing data row (135)
#        down onto the 4 new rows so the new cells pick up the same
#        cellXfs (bold/bordered id column style, date-time format column)
#        without creating any new style entries. ---
$ws.Range("A135").Copy() | Out-Null
$ws.Range("A136:A139").PasteSpecial(-4122) | Out-Null

$ws.Range("E135").Copy() | Out-Null
$ws.Range("E136:E139").PasteSpecial(-4122) | Out-Null

# --- 2. Fixture data (upcoming matches: no FTHG/FTAG/FTR result yet,
#        so columns H/I/J are intentionally left blank). ---

# Row 136: id 134
$ws.Cells.Item(136, 1).Value = 134
$ws.Cells.Item(136, 2).Value = 6788938
$ws.Cells.Item(136, 3).Value = "Croatia HNL"
$ws.Cells.Item(136, 4).Value = "Croatia HNL"
$ws.Cells.Item(136, 5).Value = 45380.58333333334
$ws.Cells.Item(136, 6).Value = "Slaven Belupo"
$ws.Cells.Item(136, 7).Value = "NK Osijek"
$ws.Cells.Item(136, 11).Value = 2.625
$ws.Cells.Item(136, 12).Value = 3.4
$ws.Cells.Item(136, 13).Value = 2.4
$ws.Cells.Item(136, 14).Value = 3.25
$ws.Cells.Item(136, 15).Value = 3.3
$ws.Cells.Item(136, 16).Value = 2.1
$ws.Cells.Item(136, 17).Value = 0.25
$ws.Cells.Item(136, 18).Value = 1.975
$ws.Cells.Item(136, 19).Value = 1.875
$ws.Cells.Item(136, 20).Value = 2.5
$ws.Cells.Item(136, 21).Value = 1.975
$ws.Cells.Item(136, 22).Value = 1.875
$ws.Cells.Item(136, 23).Value = 0
$ws.Cells.Item(136, 24).Value = 0
$ws.Cells.Item(136, 25).Value = 0
$ws.Cells.Item(136, 26).Value = 0
$ws.Cells.Item(136, 27).Value = 0

# Row 137: id 135
$ws.Cells.Item(137, 1).Value = 135
$ws.Cells.Item(137, 2).Value = 6788936
$ws.Cells.Item(137, 3).Value = "Croatia HNL"
$ws.Cells.Item(137, 4).Value = "Croatia HNL"
$ws.Cells.Item(137, 5).Value = 45381.45138888889
$ws.Cells.Item(137, 6).Value = "NK Varazdin"
$ws.Cells.Item(137, 7).Value = "HNK Gorica"
$ws.Cells.Item(137, 11).Value = 2.25
$ws.Cells.Item(137, 12).Value = 3.3
$ws.Cells.Item(137, 13).Value = 2.875
$ws.Cells.Item(137, 14).Value = 2.25
$ws.Cells.Item(137, 15).Value = 3.2
$ws.Cells.Item(137, 16).Value = 3
$ws.Cells.Item(137, 17).Value = -0.25
$ws.Cells.Item(137, 18).Value = 2.025
$ws.Cells.Item(137, 19).Value = 1.825
$ws.Cells.Item(137, 20).Value = 2.25
$ws.Cells.Item(137, 21).Value = 1.8
$ws.Cells.Item(137, 22).Value = 2.05
$ws.Cells.Item(137, 23).Value = 0
$ws.Cells.Item(137, 24).Value = 0
$ws.Cells.Item(137, 25).Value = 0
$ws.Cells.Item(137, 26).Value = 0
$ws.Cells.Item(137, 27).Value = 0

# Row 138: id 136
$ws.Cells.Item(138, 1).Value = 136
$ws.Cells.Item(138, 2).Value = 6769306
$ws.Cells.Item(138, 3).Value = "Croatia HNL"
$ws.Cells.Item(138, 4).Value = "Croatia HNL"
$ws.Cells.Item(138, 5).Value = 45381.54166666666
$ws.Cells.Item(138, 6).Value = "NK Rudes"
$ws.Cells.Item(138, 7).Value = "HNK Rijeka"
$ws.Cells.Item(138, 11).Value = 11
$ws.Cells.Item(138, 12).Value = 6
$ws.Cells.Item(138, 13).Value = 1.2
$ws.Cells.Item(138, 14).Value = 11
$ws.Cells.Item(138, 15).Value = 6
$ws.Cells.Item(138, 16).Value = 1.2
$ws.Cells.Item(138, 17).Value = 1.75
$ws.Cells.Item(138, 18).Value = 1.975
$ws.Cells.Item(138, 19).Value = 1.875
$ws.Cells.Item(138, 20).Value = 2.75
$ws.Cells.Item(138, 21).Value = 1.85
$ws.Cells.Item(138, 22).Value = 2
$ws.Cells.Item(138, 23).Value = 0
$ws.Cells.Item(138, 24).Value = 0
$ws.Cells.Item(138, 25).Value = 0
$ws.Cells.Item(138, 26).Value = 0
$ws.Cells.Item(138, 27).Value = 0

# Row 139: id 137
$ws.Cells.Item(139, 1).Value = 137
$ws.Cells.Item(139, 2).Value = 6788937
$ws.Cells.Item(139, 3).Value = "Croatia HNL"
$ws.Cells.Item(139, 4).Value = "Croatia HNL"
$ws.Cells.Item(139, 5).Value = 45381.64583333334
$ws.Cells.Item(139, 6).Value = "Hajduk Split"
$ws.Cells.Item(139, 7).Value = "Dinamo Zagreb"
$ws.Cells.Item(139, 11).Value = 2.4
$ws.Cells.Item(139, 12).Value = 3.4
$ws.Cells.Item(139, 13).Value = 2.625
$ws.Cells.Item(139, 14).Value = 2.375
$ws.Cells.Item(139, 15).Value = 3.3
$ws.Cells.Item(139, 16).Value = 2.75
$ws.Cells.Item(139, 17).Value = 0
$ws.Cells.Item(139, 18).Value = 1.775
$ws.Cells.Item(139, 19).Value = 2.1
$ws.Cells.Item(139, 20).Value = 2.25
$ws.Cells.Item(139, 21).Value = 1.825
$ws.Cells.Item(139, 22).Value = 2.025
$ws.Cells.Item(139, 23).Value = 0
$ws.Cells.Item(139, 24).Value = 0
$ws.Cells.Item(139, 25).Value = 0
$ws.Cells.Item(139, 26).Value = 0
$ws.Cells.Item(139, 27).Value = 0
